# Weekly update: insert a new weekly record as row 418 (pushing the
# existing rows 418-466 down to 419-467) on the "Hortaliza, Terminal La
# Palmera de La Serena - Zanahoria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 418; existing rows 418-466 shift down to 419-467.
$ws.Rows("418").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A418").Value = 8
$ws.Range("B418").Value = "Terminal La Palmera de La Serena"
$ws.Range("C418").Value = "Coquimbo"
$ws.Range("D418").Value = 44946
$ws.Range("E418").Value = 4
$ws.Range("F418").Value = 100114013
$ws.Range("G418").Value = "Zanahoria"
$ws.Range("H418").Value = "Sin especificar"
$ws.Range("I418").Value = "Primera"
$ws.Range("J418").Value = 600
$ws.Range("K418").Value = 5000
$ws.Range("L418").Value = 6000
$ws.Range("M418").Value = 5500
$ws.Range("N418").Value = "$/saco 20 kilos"
$ws.Range("O418").Value = "Provincia del Elquí"
$ws.Range("P418").Value = 275
$ws.Range("Q418").Value = 20
$ws.Range("R418").Value = "Hortaliza"
